# Apply the TODO list edit described by the commit:
#  - minimize the workbook window
#  - remove several completed/obsolete TODO items
#  - add two new TODO items ("call batch file" / "partial encryption")
#  - re-flow the remaining rows to close the resulting gaps

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Minimize the workbook window (workbookView minimized="1")
$excel.ActiveWindow.WindowState = -4140

# --- Row 2 ---------------------------------------------------------------
$ws.Range("B2").Value = "7) הודעה ללקוח -  ביטול רישום ל-startup processes"
$ws.Range("B2").HorizontalAlignment = -4152

# --- Row 3 ---------------------------------------------------------------
$ws.Range("B3").Value = "12) WORD  הנדסת אנוש"
$ws.Range("B3").HorizontalAlignment = -4152
$ws.Range("C3").Value = "להוסיף קריאה לbatch file."

# --- Row 4 (unchanged text, already styled) -------------------------------
$ws.Range("B4").Value = "לסגור ווינדוס דפנדר"

# --- Row 5 becomes blank (still styled) -----------------------------------
$ws.Range("B5").ClearContents()

# --- Row 6 (new row, same content that used to be row 9) ------------------
$ws.Range("B6").Value = "8) OVERALL TEST על מכונה וירטואלית"
$ws.Range("B6").HorizontalAlignment = -4152
$ws.Range("C6").Value = "כמה אנטיוירוסים"
$ws.Range("D6").Value = "בדיקת stratup"
$ws.Range("E6").Value = "בדיקת מאקרו"
$ws.Range("F6").Value = "6) BACKUP CORRUPTION"
$ws.Range("F6").HorizontalAlignment = -4152

# --- Old row 8 (now empty - content moved up to row 2) --------------------
$ws.Range("B8").Clear()

# --- Old row 9 (now empty - content moved up to row 6) --------------------
$ws.Range("B9:F9").Clear()

# --- Row 9 (new content, used to be row 13) --------------------------------
$ws.Range("B9").Value = "13) סיומת קובץ 5 תווים אקריים/ קבצי PDF עם MAGIC"
$ws.Range("B9").HorizontalAlignment = -4152

# --- Row 10 (new content, used to be row 17) -------------------------------
$ws.Range("B10").Value = "שינוי סדר איטרציה"

# --- Row 11 (new content, used to be row 18) -------------------------------
$ws.Range("B11").Value = "16) מיון קבצים נוסף לפי גודל"

# --- Row 12 becomes blank (still styled) -----------------------------------
$ws.Range("B12").ClearContents()

# --- Row 13 (new task) ------------------------------------------------------
$ws.Range("B13").Value = "הצפנה חלקית"

# --- Row 15 becomes blank (still styled) ------------------------------------
$ws.Range("B15").ClearContents()
$ws.Range("C15").ClearContents()

# --- Row 16 becomes blank (still styled) ------------------------------------
$ws.Range("B16").ClearContents()

# --- Rows 17 & 18 no longer exist -------------------------------------------
$ws.Rows("17:18").Delete()

# Move the active selection to B14 (matches the saved selection in the file)
$ws.Range("B14").Select() | Out-Null
